# Auto-generated Excel COM-interop script
# Updates scheduled-runner market/profit snapshot values across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3453.4614
$ws.Range("J38").Value = 4787.0835
$ws.Range("L38").Value = 14361.2505
$ws.Range("N38").Value = -15105.2505
$ws.Range("H100").Value = 3536.4546
$ws.Range("J100").Value = 4259.3335
$ws.Range("L100").Value = 4259.3335
$ws.Range("N100").Value = -5341.3335
$ws.Range("H133").Value = 139000
$ws.Range("J133").Value = 139000
$ws.Range("L133").Value = 139000
$ws.Range("N133").Value = -149120
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2865.4
$ws.Range("I2").Value = 1553.375
$ws.Range("J2").Value = 5197.8887
$ws.Range("K2").Value = 1553.375
$ws.Range("L2").Value = 5197.8887
$ws.Range("M2").Value = -1440.375
$ws.Range("N2").Value = -5423.8887
$ws.Range("H45").Value = 3451.9092
$ws.Range("I45").Value = 2555.75
$ws.Range("K45").Value = 2555.75
$ws.Range("M45").Value = -2178.75
$ws.Range("H74").Value = 52636984
$ws.Range("I74").Value = 58828950
$ws.Range("K74").Value = 58828950
$ws.Range("M74").Value = -58828076
$ws.Range("H77").Value = 52636984
$ws.Range("I77").Value = 58828950
$ws.Range("K77").Value = 294144750
$ws.Range("M77").Value = -294140382
$ws.Range("H81").Value = 59000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 59000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H116").Value = 2865.4
$ws.Range("I116").Value = 1553.375
$ws.Range("J116").Value = 5197.8887
$ws.Range("K116").Value = 1553.375
$ws.Range("L116").Value = 5197.8887
$ws.Range("M116").Value = 740.625
$ws.Range("N116").Value = -9785.8887
$ws.Range("H122").Value = 9525931
$ws.Range("I122").Value = 1919.8636
$ws.Range("J122").Value = 25643488
$ws.Range("K122").Value = 5759.5908
$ws.Range("L122").Value = 76930464
$ws.Range("M122").Value = -3309.5908
$ws.Range("N122").Value = -76935364
$ws.Range("H124").Value = 32248.75
$ws.Range("J124").Value = 32248.75
$ws.Range("L124").Value = 32248.75
$ws.Range("N124").Value = -42068.75
$ws.Range("H129").Value = 68749.25
$ws.Range("J129").Value = 68749.25
$ws.Range("L129").Value = 68749.25
$ws.Range("N129").Value = -78749.25
$ws.Range("H132").Value = 35776344
$ws.Range("I132").Value = 14619.318
$ws.Range("K132").Value = 43857.954
$ws.Range("M132").Value = -41327.954
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2865.4
$ws.Range("I3").Value = 1553.375
$ws.Range("J3").Value = 5197.8887
$ws.Range("K3").Value = 1553.375
$ws.Range("L3").Value = 5197.8887
$ws.Range("M3").Value = -1439.375
$ws.Range("N3").Value = -5425.8887
$ws.Range("H86").Value = 55444.777
$ws.Range("I86").Value = 41358
$ws.Range("K86").Value = 41358
$ws.Range("M86").Value = -40235
$ws.Range("H89").Value = 55444.777
$ws.Range("I89").Value = 41358
$ws.Range("K89").Value = 206790
$ws.Range("M89").Value = -201174
$ws.Range("H99").Value = 5572.5713
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H107").Value = 4750
$ws.Range("I107").Value = 4750
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4750
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -2830
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1871.2858
$ws.Range("I16").Value = 1871.2858
$ws.Range("K16").Value = 1871.2858
$ws.Range("M16").Value = -1584.2858
$ws.Range("H31").Value = 8932261
$ws.Range("I31").Value = 2930.125
$ws.Range("J31").Value = 20838036
$ws.Range("K31").Value = 2930.125
$ws.Range("L31").Value = 20838036
$ws.Range("M31").Value = -2635.125
$ws.Range("N31").Value = -20838626
$ws.Range("H34").Value = 8932261
$ws.Range("I34").Value = 2930.125
$ws.Range("J34").Value = 20838036
$ws.Range("K34").Value = 2930.125
$ws.Range("L34").Value = 20838036
$ws.Range("M34").Value = -2728.125
$ws.Range("N34").Value = -20838440
$ws.Range("H99").Value = 4215.3335
$ws.Range("J99").Value = 2896.8333
$ws.Range("L99").Value = 2896.8333
$ws.Range("N99").Value = -5892.8333
$ws.Range("H105").Value = 2304.5
$ws.Range("I105").Value = 2205.4
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 2205.4
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -458.4000000000001
$ws.Range("N105").Value = -6294
$ws.Range("H113").Value = 1871.2858
$ws.Range("I113").Value = 1871.2858
$ws.Range("K113").Value = 1871.2858
$ws.Range("M113").Value = 298.7141999999999
$ws.Range("H126").Value = 4215.3335
$ws.Range("J126").Value = 2896.8333
$ws.Range("L126").Value = 8690.499899999999
$ws.Range("N126").Value = -13630.4999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48758956
$ws.Range("I4").Value = 74460830
$ws.Range("J4").Value = 27592706
$ws.Range("K4").Value = 223382490
$ws.Range("L4").Value = 82778118
$ws.Range("M4").Value = -223382378
$ws.Range("N4").Value = -82778342
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11856.857
$ws.Range("I102").Value = 8999.5
$ws.Range("K102").Value = 8999.5
$ws.Range("M102").Value = -7377.5
$ws.Range("H132").Value = 2964.6
$ws.Range("I132").Value = 2955.75
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8867.25
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6337.25
$ws.Range("N132").Value = -14060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5016.636
$ws.Range("I68").Value = 3909.7778
$ws.Range("K68").Value = 3909.7778
$ws.Range("M68").Value = -3160.7778
$ws.Range("H71").Value = 5016.636
$ws.Range("I71").Value = 3909.7778
$ws.Range("K71").Value = 19548.889
$ws.Range("M71").Value = -15804.889
$ws.Range("H132").Value = 3433
$ws.Range("I132").Value = 3289.1072
$ws.Range("K132").Value = 9867.321599999999
$ws.Range("M132").Value = -7337.321599999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2823.9167
$ws.Range("I81").Value = 1482.1666
$ws.Range("K81").Value = 2964.3332
$ws.Range("M81").Value = -1903.3332
$ws.Range("H84").Value = 2823.9167
$ws.Range("I84").Value = 1482.1666
$ws.Range("K84").Value = 14821.666
$ws.Range("M84").Value = -9517.666000000001
$ws.Range("H124").Value = 2526500
$ws.Range("J124").Value = 2526500
$ws.Range("L124").Value = 2526500
$ws.Range("N124").Value = -2536320
$ws.Range("H125").Value = 65368.168
$ws.Range("J125").Value = 65368.168
$ws.Range("L125").Value = 65368.168
$ws.Range("N125").Value = -75208.16800000001
$ws.Range("H132").Value = 4338.778
$ws.Range("I132").Value = 3721.3572
$ws.Range("K132").Value = 11164.0716
$ws.Range("M132").Value = -9105.3335
